# Repartition des membres dans chaque equipe : 7/7/4 -> 6/6/6
# On prend la liste a plat des membres (colonne A, puis B, puis C, dans
# l'ordre des lignes) et on la repartage en 3 groupes de 6, en abandonnant
# les 3 derniers membres (F. Alonso, M. Singh, D. Wade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colonnes utilisees pour les equipes : C (A), D (B), E (C)
$teamCols = @(3, 4, 5)

# Lire tous les membres actuels (lignes 2 a 8) dans l'ordre C, D, E
$members = @()
foreach ($col in $teamCols) {
    for ($row = 2; $row -le 8; $row++) {
        $val = $ws.Cells.Item($row, $col).Value2
        if ($val -ne $null -and $val -ne "") {
            $members += $val
        }
    }
}

# Ne garder que les 18 premiers (6 par equipe), les 3 derniers sont retires
$members = $members[0..17]

# Effacer les anciennes valeurs (lignes 2 a 8, colonnes C a E)
$ws.Range("C2:E8").ClearContents()

# Re-ecrire les membres repartis en 6/6/6
$idx = 0
foreach ($col in $teamCols) {
    for ($row = 2; $row -le 7; $row++) {
        $ws.Cells.Item($row, $col).Value = $members[$idx]
        $idx++
    }
}
